# ProjectPlanner.xlsx update
# - Updates G8 and G14 percentages
# - Reworks the "Testing" sub-section (rows 19-28): adds a new "Algorithm" activity,
#   shifts "Remote" and its children down by one row, updates owners/periods/percentages
#   for "Build GUI" and "Testing" rows, and fixes the dependent AVERAGE() formulas.
# - Adjusts the sheet view's visible top row / selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- simple % complete tweaks -------------------------------------------------
$ws.Range("G8").Value2  = 0.5
$ws.Range("G14").Value2 = 0.9

# --- "Testing" header (row 19) ------------------------------------------------
# Owner changes from Ruben to Håkon, and the roll-up average now also covers the
# new row 22 (Algorithm).
$ws.Range("A19").Value2 = "Håkon"
$ws.Range("G19").Formula = "=AVERAGE(G20:G22)"

# --- new row 22: "Algorithm" (child of Testing) -------------------------------
# Pick up the look of an existing "child" row (e.g. row 21) so the new row matches
# the rest of the sub-items, then fill in its own content.
$ws.Range("A21:G21").Copy() | Out-Null
$ws.Range("A22:G22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A22").Value2 = ""
$ws.Range("B22").Value2 = "Algorithm"
$ws.Range("C22").Value2 = 43
$ws.Range("D22").Value2 = 1
$ws.Range("E22").Value2 = 43
$ws.Range("F22").Value2 = ""
$ws.Range("G22").Value2 = 0.3

# --- row 23: becomes the "Remote" header (was row 22 content) ----------------
$ws.Range("A19:G19").Copy() | Out-Null
$ws.Range("A23:G23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A23").Value2 = ""
$ws.Range("B23").Value2 = "Remote"
$ws.Range("C23").Value2 = 43
$ws.Range("D23").Value2 = 2
$ws.Range("E23").Value2 = ""
$ws.Range("F23").Value2 = ""
$ws.Range("G23").Formula = "=AVERAGE(G24:G26)"

# --- row 24: "Logging" (child) ------------------------------------------------
$ws.Range("B24").Value2 = "Logging"
$ws.Range("C24").Value2 = 43
$ws.Range("D24").Value2 = 1
$ws.Range("E24").Value2 = ""
$ws.Range("G24").Value2 = 0

# --- row 25: "Controlling" (child) -------------------------------------------
$ws.Range("B25").Value2 = "Controlling"
$ws.Range("C25").Value2 = 43
$ws.Range("D25").Value2 = 1
$ws.Range("E25").Value2 = ""
$ws.Range("G25").Value2 = 0

# --- row 26: "Monitoring" (child) --------------------------------------------
$ws.Range("A26").Value2 = ""
$ws.Range("B26").Value2 = "Monitoring"
$ws.Range("C26").Value2 = 44
$ws.Range("D26").Value2 = 1
$ws.Range("E26").Value2 = ""
$ws.Range("G26").Value2 = 0

# --- row 27: "Build GUI" header (Ruben) --------------------------------------
$ws.Range("A27").Value2 = "Ruben"
$ws.Range("B27").Value2 = "Build GUI"
$ws.Range("C27").Value2 = 45
$ws.Range("D27").Value2 = 1
$ws.Range("E27").Value2 = 43
$ws.Range("F27").Value2 = ""
$ws.Range("G27").Value2 = 0.4

# --- row 28: "Testing" header (ALL) ------------------------------------------
$ws.Range("A27:G27").Copy() | Out-Null
$ws.Range("A28:G28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A28").Value2 = "ALL"
$ws.Range("B28").Value2 = "Testing"
$ws.Range("C28").Value2 = 47
$ws.Range("D28").Value2 = 2
$ws.Range("E28").Value2 = ""
$ws.Range("F28").Value2 = ""
$ws.Range("G28").Value2 = 0

# --- sheet view: scroll & selection -------------------------------------------
$ws.Range("A12").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("G28").Select() | Out-Null

Write-Host "ProjectPlanner.xlsx Testing section updated."
